$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header capitalization: Cargo_type -> Cargo_Type
$ws.Range("A1").Value = "Cargo_Type"

# Update Cargo type data values from GEN -> GENERAL
$ws.Range("A2:A9").Value = "GENERAL"

# Add a new "Tax" column header, styled like the other secondary header cells (D1/E1/O1/P1)
$ws.Range("D1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q1").Value = "Tax"

# Row heights picked up by Excel's autofit after the edits
$ws.Rows.Item(1).RowHeight = 25.5
$ws.Range("A2:A9").RowHeight = 15.75

$null = $ws.Range("I21").Select()
